$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores everything as literal text (even plain-looking
# numbers like "1.00"), matching the source feed. Excel auto-detects
# numeric-looking input on a General-formatted cell and converts it to a
# real number, which would corrupt values such as "1.00" -> 1 or
# "26.90" -> 26.9. Pre-format those specific cells as Text so the literal
# string is preserved, exactly like the original cells already are.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.373.89'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '3.677.76'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '684.73'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").Value = '158.97'
$ws.Range("E6").Value = '  -2.32%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -1.35%  '
$ws.Range("E9").Value = '  -2.32%  '
$ws.Range("D10").Value = '7.03'
$ws.Range("E10").Value = '  -3.32%  '
$ws.Range("D11").Value = '0.435'
$ws.Range("E11").Value = '  -3.34%  '
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("D13").Value = '4.299.45'
$ws.Range("E13").Value = '  -0.46%  '
$ws.Range("D14").Value = '32.25'
$ws.Range("E14").Value = '  -3.99%  '
$ws.Range("D15").Value = '3.670.36'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").Value = '69.389.00'
$ws.Range("E17").Value = '  +2.03%  '
$ws.Range("D18").Value = '15.83'
$ws.Range("E18").Value = '  -3.27%  '
$ws.Range("D19").Value = '6.38'
$ws.Range("E19").Value = '  -3.57%  '
$ws.Range("D20").Value = '469.74'
$ws.Range("E20").Value = '  -3.03%  '
$ws.Range("D21").Value = '9.94'
$ws.Range("E21").Value = '  +1.36%  '
$ws.Range("D22").Value = '0.649'
$ws.Range("E22").Value = '  -2.78%  '
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("D24").Value = '3.822.95'
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("D26").Value = '0.0000123'
$ws.Range("E26").Value = '  -4.38%  '
$ws.Range("D27").Value = '10.92'
$ws.Range("E27").Value = '  -5.41%  '
$ws.Range("E28").Value = '  -4.21%  '
$ws.Range("E29").Value = '  -1.80%  '
$ws.Range("D30").Value = '1.73'
$ws.Range("E30").Value = '  -6.13%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = '6.54'
$ws.Range("E32").Value = '  -4.12%  '
$ws.Range("E33").Value = '  -6.33%  '
$ws.Range("E34").Value = '  -0.83%  '
$ws.Range("D35").Value = '3.652.01'
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("E36").Value = '  -3.44%  '
$ws.Range("D37").Value = '8.16'
$ws.Range("E37").Value = '  -4.83%  '
$ws.Range("E38").Value = '  +1.39%  '
$ws.Range("D40").Value = '2.23'
$ws.Range("E40").Value = '  +2.21%  '
$ws.Range("D41").Value = '0.0898'
$ws.Range("E41").Value = '  -4.82%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").Value = '166.22'
$ws.Range("E43").Value = '  +5.43%  '
$ws.Range("E44").Value = '  -2.48%  '
$ws.Range("E45").Value = '  -1.43%  '
$ws.Range("D46").Value = '0.000283'
$ws.Range("E46").Value = '  +1.12%  '
$ws.Range("D47").Value = '2.71'
$ws.Range("E47").Value = '  -4.65%  '
$ws.Range("E48").Value = '  +3.99%  '
$ws.Range("E49").Value = '  +0.25%  '
$ws.Range("D50").Value = '27.31'
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("D51").Value = '7.78'
$ws.Range("E51").Value = '  -4.01%  '
